$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "AEDB.CEA"
$ws.Range("B2").Value = "epmajor.30days"
$ws.Range("C2").Value = "MCP1_pg_ug_2015_rank"
$ws.Range("D2").Value = 1.34745908919635
$ws.Range("E2").Value = 0.424795046279012
$ws.Range("F2").Value = 3.84763659816062
$ws.Range("G2").Value = 1.67340554500438
$ws.Range("H2").Value = 8.84681387348115
$ws.Range("I2").Value = 3.17202166315123
$ws.Range("J2").Value = 0.00151381714207261
$ws.Range("K2").Value = 1029
$ws.Range("L2").Value = 32

# Row 3
$ws.Range("A3").Value = "AEDB.CEA"
$ws.Range("B3").Value = "epmajor.30days"
$ws.Range("C3").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("D3").Value = 1.16724283964527
$ws.Range("E3").Value = 0.418059106149549
$ws.Range("F3").Value = 3.21312132360115
$ws.Range("G3").Value = 1.41601557087562
$ws.Range("H3").Value = 7.29098524940393
$ws.Range("I3").Value = 2.79205218227616
$ws.Range("J3").Value = 0.00523749089470836
$ws.Range("K3").Value = 1029
$ws.Range("L3").Value = 32

# Row 4
$ws.Range("A4").Value = "AEDB.CEA"
$ws.Range("B4").Value = "epmajor.30days"
$ws.Range("C4").Value = "MCP1_rank"
$ws.Range("D4").Value = -0.0678094579469317
$ws.Range("E4").Value = 0.579574510320059
$ws.Range("F4").Value = 0.934438506446785
$ws.Range("G4").Value = 0.300059197685348
$ws.Range("H4").Value = 2.91001018821006
$ws.Range("I4").Value = -0.116998689106402
$ws.Range("J4").Value = 0.906861092289218
$ws.Range("K4").Value = 493
$ws.Range("L4").Value = 14

# Row 5
$ws.Range("A5").Value = "AEDB.CEA"
$ws.Range("B5").Value = "epstroke.30days"
$ws.Range("C5").Value = "MCP1_pg_ug_2015_rank"
$ws.Range("D5").Value = 1.10546122275982
$ws.Range("E5").Value = 0.446587176506691
$ws.Range("F5").Value = 3.02061732503239
$ws.Range("G5").Value = 1.25878926323529
$ws.Range("H5").Value = 7.24833718460175
$ws.Range("I5").Value = 2.4753537067656
$ws.Range("J5").Value = 0.013310428009947
$ws.Range("K5").Value = 1029
$ws.Range("L5").Value = 27

# Row 6
$ws.Range("A6").Value = "AEDB.CEA"
$ws.Range("B6").Value = "epstroke.30days"
$ws.Range("C6").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("D6").Value = 0.924775855993484
$ws.Range("E6").Value = 0.443385308912267
$ws.Range("F6").Value = 2.52130306204764
$ws.Range("G6").Value = 1.05732343189811
$ws.Range("H6").Value = 6.01232218913255
$ws.Range("I6").Value = 2.08571605194179
$ws.Range("J6").Value = 0.0370043515312526
$ws.Range("K6").Value = 1029
$ws.Range("L6").Value = 27

# Row 7
$ws.Range("A7").Value = "AEDB.CEA"
$ws.Range("B7").Value = "epstroke.30days"
$ws.Range("C7").Value = "MCP1_rank"
$ws.Range("D7").Value = 0.462099695870797
$ws.Range("E7").Value = 0.622962571046731
$ws.Range("F7").Value = 1.58740355291656
$ws.Range("G7").Value = 0.468177892372908
$ws.Range("H7").Value = 5.38224909988921
$ws.Range("I7").Value = 0.741777624126526
$ws.Range("J7").Value = 0.458222079798769
$ws.Range("K7").Value = 493
$ws.Range("L7").Value = 12

# Row 8
$ws.Range("A8").Value = "AEDB.CEA"
$ws.Range("B8").Value = "epcoronary.30days"
$ws.Range("C8").Value = "MCP1_pg_ug_2015_rank"
$ws.Range("D8").Value = 1.16987216279731
$ws.Range("E8").Value = 0.874703845602647
$ws.Range("F8").Value = 3.22158077432883
$ws.Range("G8").Value = 0.580104293873301
$ws.Range("H8").Value = 17.8908909917359
$ws.Range("I8").Value = 1.33744943351803
$ws.Range("J8").Value = 0.181075977622402
$ws.Range("K8").Value = 1029
$ws.Range("L8").Value = 8

# Row 9
$ws.Range("A9").Value = "AEDB.CEA"
$ws.Range("B9").Value = "epcoronary.30days"
$ws.Range("C9").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("D9").Value = 1.98230046039678
$ws.Range("E9").Value = 1.10915521576149
$ws.Range("F9").Value = 7.25942380799215
$ws.Range("G9").Value = 0.825600926124824
$ws.Range("H9").Value = 63.8313649566759
$ws.Range("I9").Value = 1.78721646188701
$ws.Range("J9").Value = 0.073902509662866
$ws.Range("K9").Value = 1029
$ws.Range("L9").Value = 8

# Row 10
$ws.Range("A10").Value = "AEDB.CEA"
$ws.Range("B10").Value = "epcoronary.30days"
$ws.Range("C10").Value = "MCP1_rank"
$ws.Range("D10").Value = -1.12174461052695
$ws.Range("E10").Value = 1.23135982374534
$ws.Range("F10").Value = 0.325711059713317
$ws.Range("G10").Value = 0.0291526382269638
$ws.Range("H10").Value = 3.63904266892215
$ws.Range("I10").Value = -0.910980355940977
$ws.Range("J10").Value = 0.362305724141779
$ws.Range("K10").Value = 493
$ws.Range("L10").Value = 3

# Row 11
$ws.Range("A11").Value = "AEDB.CEA"
$ws.Range("B11").Value = "epcvdeath.30days"
$ws.Range("C11").Value = "MCP1_pg_ug_2015_rank"
$ws.Range("D11").Value = 54.4935784085021
$ws.Range("E11").Value = 1719.22673097905
$ws.Range("F11").Value = [double]"463724884771023005679616"
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = "#NUM!"
$ws.Range("I11").Value = 0.0316965630109006
$ws.Range("J11").Value = 0.974714035831114
$ws.Range("K11").Value = 1029
$ws.Range("L11").Value = 2

# Row 12
$ws.Range("A12").Value = "AEDB.CEA"
$ws.Range("B12").Value = "epcvdeath.30days"
$ws.Range("C12").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("D12").Value = 53.4330043812754
$ws.Range("E12").Value = 536.750098063802
$ws.Range("F12").Value = [double]"160567983869854983651328"
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = "#NUM!"
$ws.Range("I12").Value = 0.0995491283076094
$ws.Range("J12").Value = 0.920702282837156
$ws.Range("K12").Value = 1029
$ws.Range("L12").Value = 2

# Row 13
$ws.Range("A13").Value = "AEDB.CEA"
$ws.Range("B13").Value = "epcvdeath.30days"
$ws.Range("C13").Value = "MCP1_rank"
$ws.Range("D13").Value = 24.9249999162542
$ws.Range("E13").Value = 46394.8743112099
$ws.Range("F13").Value = 66802070749.6082
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = "#NUM!"
$ws.Range("I13").Value = 0.000537236069421398
$ws.Range("J13").Value = 0.999571347655322
$ws.Range("K13").Value = 493
$ws.Range("L13").Value = 1
